$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed price/volume snapshot for the cryptos list (GitHub Actions bot).
# Helper: some "Price" strings are plain decimals (e.g. "601.71") that Excel
# would otherwise auto-convert to a Double when assigned via .Value. Prefixing
# with an apostrophe forces text storage, matching the source sheet, which keeps
# every Price cell a string (it also uses "12.345.67"-style grouped numbers that
# are never valid numerics anyway).
function Set-CellText($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range('D2').Value = '69.996.32'
$ws.Range('E2').Value = '  +2.09%  '
$ws.Range('D3').Value = '3.506.33'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-CellText 'D5' '601.71'
$ws.Range('E5').Value = '  +2.61%  '
Set-CellText 'D6' '171.44'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').Value = '3.500.74'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  +7.26%  '
Set-CellText 'D12' '0.579'
$ws.Range('E12').Value = '  +1.39%  '
Set-CellText 'D13' '45.99'
$ws.Range('E13').Value = '  -0.91%  '
Set-CellText 'D14' '0.0000274'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '4.063.16'
$ws.Range('E15').Value = '  +0.75%  '
Set-CellText 'D16' '8.24'
$ws.Range('E16').Value = '  -0.70%  '
Set-CellText 'D17' '603.38'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.032.32'
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.490.48'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  +0.97%  '
Set-CellText 'D21' '17.08'
$ws.Range('E21').Value = '  -0.67%  '
Set-CellText 'D22' '0.866'
$ws.Range('E22').Value = '  -0.55%  '
Set-CellText 'D23' '9.14'
$ws.Range('E23').Value = '  -17.03%  '
Set-CellText 'D24' '15.46'
$ws.Range('E24').Value = '  -1.52%  '
Set-CellText 'D25' '95.27'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('E26').Value = '  -1.70%  '
$ws.Range('E27').Value = '  -0.02%  '
Set-CellText 'D28' '2.56'
$ws.Range('E28').Value = '  -1.02%  '
Set-CellText 'D29' '33.92'
$ws.Range('E29').Value = '  +3.70%  '
Set-CellText 'D30' '8.94'
$ws.Range('E30').Value = '  -1.83%  '
Set-CellText 'D31' '692.03'
$ws.Range('E31').Value = '  +20.03%  '
$ws.Range('E32').Value = '  -2.18%  '
Set-CellText 'D33' '8.06'
$ws.Range('E33').Value = '  -3.76%  '
Set-CellText 'D34' '6.87'
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText 'D36' '3.55'
$ws.Range('E36').Value = '  +2.31%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText 'D37' '0.0993'
$ws.Range('E37').Value = '  -1.41%  '
Set-CellText 'D38' '10.65'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('E39').Value = '  +7.93%  '
Set-CellText 'D40' '56.54'
$ws.Range('E40').Value = '  -0.27%  '
Set-CellText 'D41' '0.999'
$ws.Range('E41').Value = '  -0.08%  '
Set-CellText 'D42' '0.142'
$ws.Range('E42').Value = '  +4.47%  '
$ws.Range('D43').Value = '3.319.19'
$ws.Range('E43').Value = '  -2.41%  '
Set-CellText 'D44' '0.312'
$ws.Range('E44').Value = '  -3.20%  '
Set-CellText 'D45' '2.92'
$ws.Range('E45').Value = '  +4.45%  '
Set-CellText 'D46' '32.10'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '0.0₃0686'
$ws.Range('E47').Value = '  -0.33%  '
Set-CellText 'D48' '2.54'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +0.87%  '
Set-CellText 'D50' '133.04'
$ws.Range('E50').Value = '  +0.53%  '
